$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.861638307571411
$ws.Range("B1").Value = 5.746891021728516
$ws.Range("C1").Value = 4.767726898193359
$ws.Range("D1").Value = 5.553970813751221
$ws.Range("E1").Value = 4.398074150085449
